$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force H1179:J1199 to be stored as text so values like "$ 50" and "5"
# are written as literal strings instead of being auto-converted to numbers.
$ws.Range("H1179:J1199").NumberFormat = "@"

# Row 1179
$ws.Range("A1179").Value = '21CRB01268'
$ws.Range("B1179").Value = 'Hemmeter'
$ws.Range("C1179").Value = 'Possession Drug Paraphernalia'
$ws.Range("D1179").Value = '2925.14(C)'
$ws.Range("E1179").Value = 'M4'
$ws.Range("F1179").Value = 'No Contest'
$ws.Range("G1179").Value = 'Guilty'
$ws.Range("H1179").Value = '$ 50'
$ws.Range("I1179").Value = '$ 0'
$ws.Range("J1179").Value = '5'
$ws.Range("K1179").Value = 'None'

# Row 1180
$ws.Range("A1180").Value = '21TRD09246'
$ws.Range("B1180").Value = 'Hemmeter'
$ws.Range("C1180").Value = '1st Speed In 1 Yr >70mph'
$ws.Range("D1180").Value = '4511.21D4'
$ws.Range("E1180").Value = 'No Data'
$ws.Range("F1180").Value = 'No Contest'
$ws.Range("G1180").Value = 'Guilty'
$ws.Range("H1180").Value = '$ 50'
$ws.Range("I1180").Value = '$ 0'
$ws.Range("J1180").Value = '5'
$ws.Range("K1180").Value = 'None'

# Row 1181
$ws.Range("A1181").Value = '21TRD09246'
$ws.Range("B1181").Value = 'Hemmeter'
$ws.Range("C1181").Value = '1st Speed In 1 Yr >70mph'
$ws.Range("D1181").Value = '4511.21D4'
$ws.Range("E1181").Value = 'No Data'
$ws.Range("F1181").Value = 'No Contest'
$ws.Range("G1181").Value = 'Guilty'
$ws.Range("H1181").Value = '$ 50'
$ws.Range("I1181").Value = '$ 0'
$ws.Range("J1181").Value = '5'
$ws.Range("K1181").Value = 'None'

# Row 1182
$ws.Range("A1182").Value = '21TRD09246'
$ws.Range("B1182").Value = 'Hemmeter'
$ws.Range("C1182").Value = '1st Speed In 1 Yr >70mph'
$ws.Range("D1182").Value = '4511.21D4'
$ws.Range("E1182").Value = 'No Data'
$ws.Range("F1182").Value = 'No Contest'
$ws.Range("G1182").Value = 'Guilty'
$ws.Range("H1182").Value = '$ 50'
$ws.Range("I1182").Value = '$ 0'
$ws.Range("J1182").Value = '5'
$ws.Range("K1182").Value = 'None'

# Row 1183
$ws.Range("A1183").Value = '21CRB01291'
$ws.Range("B1183").Value = 'Hemmeter'
$ws.Range("C1183").Value = 'Permission Req''d To Use Licensed Dock'
$ws.Range("D1183").Value = '1501:46-12-04'
$ws.Range("E1183").Value = 'MM'
$ws.Range("F1183").Value = 'No Contest'
$ws.Range("G1183").Value = 'Guilty'
$ws.Range("H1183").Value = '$ 50'
$ws.Range("I1183").Value = '$ 0'
$ws.Range("J1183").Value = '5'
$ws.Range("K1183").Value = 'None'

# Row 1184
$ws.Range("A1184").Value = '21CRB01291'
$ws.Range("B1184").Value = 'Hemmeter'
$ws.Range("C1184").Value = 'Permission Req''d To Use Licensed Dock'
$ws.Range("D1184").Value = '1501:46-12-04'
$ws.Range("E1184").Value = 'MM'
$ws.Range("F1184").Value = 'No Contest'
$ws.Range("G1184").Value = 'Guilty'
$ws.Range("H1184").Value = '$ 50'
$ws.Range("I1184").Value = '$ 0'
$ws.Range("J1184").Value = '5'
$ws.Range("K1184").Value = 'None'

# Row 1185
$ws.Range("A1185").Value = '21CRB01291'
$ws.Range("B1185").Value = 'Hemmeter'
$ws.Range("C1185").Value = 'Permission Req''d To Use Licensed Dock'
$ws.Range("D1185").Value = '1501:46-12-04'
$ws.Range("E1185").Value = 'MM'
$ws.Range("F1185").Value = 'No Contest'
$ws.Range("G1185").Value = 'Guilty'
$ws.Range("H1185").Value = '$ 50'
$ws.Range("I1185").Value = '$ 0'
$ws.Range("J1185").Value = '5'
$ws.Range("K1185").Value = 'None'

# Row 1186
$ws.Range("A1186").Value = '21CRB01387'
$ws.Range("B1186").Value = 'Hemmeter'
$ws.Range("C1186").Value = 'Sexual Imposition M1'
$ws.Range("D1186").Value = '2907.06(A)(1)'
$ws.Range("E1186").Value = 'M1'
$ws.Range("F1186").Value = 'No Contest'
$ws.Range("G1186").Value = 'Guilty'
$ws.Range("H1186").Value = '$ 50'
$ws.Range("I1186").Value = '$ 0'
$ws.Range("J1186").Value = '5'
$ws.Range("K1186").Value = 'None'

# Row 1187
$ws.Range("A1187").Value = '21CRB01387'
$ws.Range("B1187").Value = 'Hemmeter'
$ws.Range("C1187").Value = 'Sexual Imposition M1'
$ws.Range("D1187").Value = '2907.06(A)(1)'
$ws.Range("E1187").Value = 'M1'
$ws.Range("F1187").Value = 'No Contest'
$ws.Range("G1187").Value = 'Guilty'
$ws.Range("H1187").Value = '$ 50'
$ws.Range("I1187").Value = '$ 0'
$ws.Range("J1187").Value = '5'
$ws.Range("K1187").Value = 'None'

# Row 1188
$ws.Range("A1188").Value = '21CRB01437'
$ws.Range("B1188").Value = 'Hemmeter'
$ws.Range("C1188").Value = 'Possession Of Marihuana'
$ws.Range("D1188").Value = '2925.11C3'
$ws.Range("E1188").Value = 'MM'
$ws.Range("F1188").Value = 'No Contest'
$ws.Range("G1188").Value = 'Guilty'
$ws.Range("H1188").Value = '$ 50'
$ws.Range("I1188").Value = '$ 0'
$ws.Range("J1188").Value = '5'
$ws.Range("K1188").Value = 'None'

# Row 1189
$ws.Range("A1189").Value = '21CRB01437'
$ws.Range("B1189").Value = 'Hemmeter'
$ws.Range("C1189").Value = 'Possession Of Marihuana'
$ws.Range("D1189").Value = '2925.11C3'
$ws.Range("E1189").Value = 'MM'
$ws.Range("F1189").Value = 'No Contest'
$ws.Range("G1189").Value = 'Guilty'
$ws.Range("H1189").Value = '$ 50'
$ws.Range("I1189").Value = '$ 0'
$ws.Range("J1189").Value = '5'
$ws.Range("K1189").Value = 'None'

# Row 1190
$ws.Range("A1190").Value = '21CRB00626'
$ws.Range("B1190").Value = 'Hemmeter'
$ws.Range("C1190").Value = 'Criminal Mischief'
$ws.Range("D1190").Value = '2909.07(A)(1)'
$ws.Range("E1190").Value = 'M3'
$ws.Range("F1190").Value = 'No Contest'
$ws.Range("G1190").Value = 'Guilty'
$ws.Range("H1190").Value = '$ 50'
$ws.Range("I1190").Value = '$ 0'
$ws.Range("J1190").Value = '5'
$ws.Range("K1190").Value = 'None'

# Row 1191
$ws.Range("A1191").Value = '21CRB00626'
$ws.Range("B1191").Value = 'Hemmeter'
$ws.Range("C1191").Value = 'Assault - M1'
$ws.Range("D1191").Value = '2903.13(A)'
$ws.Range("E1191").Value = 'M1'
$ws.Range("F1191").Value = 'No Contest'
$ws.Range("G1191").Value = 'Guilty'
$ws.Range("H1191").Value = '$ 0'
$ws.Range("I1191").Value = '$ 0'
$ws.Range("J1191").Value = 'None'
$ws.Range("K1191").Value = 'None'

# Row 1192
$ws.Range("A1192").Value = '21CRB00626'
$ws.Range("B1192").Value = 'Hemmeter'
$ws.Range("C1192").Value = 'Disorderly Conduct'
$ws.Range("D1192").Value = '2917.11A1'
$ws.Range("E1192").Value = 'MM'
$ws.Range("F1192").Value = 'No Contest'
$ws.Range("G1192").Value = 'Guilty'
$ws.Range("H1192").Value = '$ 0'
$ws.Range("I1192").Value = '$ 0'
$ws.Range("J1192").Value = 'None'
$ws.Range("K1192").Value = 'None'

# Row 1193
$ws.Range("A1193").Value = '21CRB01437'
$ws.Range("B1193").Value = 'Hemmeter'
$ws.Range("C1193").Value = 'Possession Of Marihuana'
$ws.Range("D1193").Value = '2925.11C3'
$ws.Range("E1193").Value = 'MM'
$ws.Range("F1193").Value = 'No Contest'
$ws.Range("G1193").Value = 'Guilty'
$ws.Range("H1193").Value = '$ 0'
$ws.Range("I1193").Value = '$ 0'
$ws.Range("J1193").Value = '5'
$ws.Range("K1193").Value = 'None'

# Row 1194
$ws.Range("A1194").Value = '21CRB01437'
$ws.Range("B1194").Value = 'Hemmeter'
$ws.Range("C1194").Value = 'Possession Of Marihuana'
$ws.Range("D1194").Value = '2925.11C3'
$ws.Range("E1194").Value = 'MM'
$ws.Range("F1194").Value = 'No Contest'
$ws.Range("G1194").Value = 'Guilty'
$ws.Range("H1194").Value = '$ 50'
$ws.Range("I1194").Value = '$ 0'
$ws.Range("J1194").Value = '5'
$ws.Range("K1194").Value = 'None'

# Row 1195
$ws.Range("A1195").Value = '21CRB01437'
$ws.Range("B1195").Value = 'Hemmeter'
$ws.Range("C1195").Value = 'Possession Of Marihuana'
$ws.Range("D1195").Value = '2925.11C3'
$ws.Range("E1195").Value = 'MM'
$ws.Range("F1195").Value = 'No Contest'
$ws.Range("G1195").Value = 'Guilty'
$ws.Range("H1195").Value = '$ 50'
$ws.Range("I1195").Value = '$ 0'
$ws.Range("J1195").Value = '5'
$ws.Range("K1195").Value = 'None'

# Row 1196
$ws.Range("A1196").Value = '21CRB01437'
$ws.Range("B1196").Value = 'Hemmeter'
$ws.Range("C1196").Value = 'Possession Of Marihuana'
$ws.Range("D1196").Value = '2925.11C3'
$ws.Range("E1196").Value = 'MM'
$ws.Range("F1196").Value = 'No Contest'
$ws.Range("G1196").Value = 'Guilty'
$ws.Range("H1196").Value = '$ 0'
$ws.Range("I1196").Value = '$ 0'
$ws.Range("J1196").Value = '5'
$ws.Range("K1196").Value = 'None'

# Row 1197
$ws.Range("A1197").Value = '21CRB01437'
$ws.Range("B1197").Value = 'Hemmeter'
$ws.Range("C1197").Value = 'Possession Of Marihuana'
$ws.Range("D1197").Value = '2925.11C3'
$ws.Range("E1197").Value = 'MM'
$ws.Range("F1197").Value = 'No Contest'
$ws.Range("G1197").Value = 'Guilty'
$ws.Range("H1197").Value = '$ 33'
$ws.Range("I1197").Value = '$ 0'
$ws.Range("J1197").Value = '3'
$ws.Range("K1197").Value = 'None'

# Row 1198
$ws.Range("A1198").Value = '21CRB01437'
$ws.Range("B1198").Value = 'Hemmeter'
$ws.Range("C1198").Value = 'Possession Of Marihuana'
$ws.Range("D1198").Value = '2925.11C3'
$ws.Range("E1198").Value = 'MM'
$ws.Range("F1198").Value = 'No Contest'
$ws.Range("G1198").Value = 'Guilty'
$ws.Range("H1198").Value = '$ 33'
$ws.Range("I1198").Value = '$ 0'
$ws.Range("J1198").Value = '3'
$ws.Range("K1198").Value = 'None'

# Row 1199
$ws.Range("A1199").Value = '21CRB01437'
$ws.Range("B1199").Value = 'Hemmeter'
$ws.Range("C1199").Value = 'Possession Of Marihuana'
$ws.Range("D1199").Value = '2925.11C3'
$ws.Range("E1199").Value = 'MM'
$ws.Range("F1199").Value = 'No Contest'
$ws.Range("G1199").Value = 'Guilty'
$ws.Range("H1199").Value = '$ 4'
$ws.Range("I1199").Value = '$ 0'
$ws.Range("J1199").Value = '2'
$ws.Range("K1199").Value = 'None'
